$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 44 is Idaho; fill in the previously-missing data now that the run succeeded.
$ws.Range("B44").NumberFormat = "YYYY-MM-DD"
$ws.Range("B44").Value = (Get-Date -Year 2020 -Month 7 -Day 22 -Hour 0 -Minute 0 -Second 0)

$ws.Range("C44").Value = 16322
$ws.Range("D44").Value = 131
$ws.Range("E44").Value = 208
$ws.Range("F44").Value = 2
$ws.Range("G44").Value = 1.27
$ws.Range("H44").Value = 1.53

$ws.Range("J44").Value = $true

$ws.Range("O44").Value = "Success!"
